# Atualizando o arquivo XLSX
# Updates odds values on Sheet1 to match the latest FlashScore scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Los Andes vs Alvarado)
$ws.Range("N3").Value = 2.87
$ws.Range("O3").Value = 1.37

# Row 4 (CA Estudiantes vs Temperley)
$ws.Range("O4").Value = 1.3
$ws.Range("R4").Value = 2.62
$ws.Range("S4").Value = 1.41

# Row 5 (Ekenas vs Lahti)
$ws.Range("L5").Value = 1.29
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 1.85
$ws.Range("O5").Value = 1.75
$ws.Range("P5").Value = 1.4
$ws.Range("Q5").Value = 2.52
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.98
$ws.Range("T5").Value = 9
$ws.Range("U5").Value = 14
$ws.Range("Z5").Value = 9.75
$ws.Range("AA5").Value = 6.2
$ws.Range("AB5").Value = 13
$ws.Range("AD5").Value = 8.25
$ws.Range("AE5").Value = 12.5
$ws.Range("AF5").Value = 9.5
$ws.Range("AG5").Value = 27
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 29
$ws.Range("AJ5").Value = 450

# Row 6 (JaPS vs KaPa)
$ws.Range("H6").Value = 4.25
$ws.Range("I6").Value = 4.15
$ws.Range("S6").Value = 2.95
$ws.Range("AH6").Value = 32

# Row 10 (Guairena vs Pastoreo)
$ws.Range("G10").Value = 2.25
$ws.Range("H10").Value = 3
$ws.Range("L10").Value = 1.47
$ws.Range("M10").Value = 2.35
$ws.Range("Q10").Value = 2.18
$ws.Range("U10").Value = 9.5
$ws.Range("V10").Value = 9.75
$ws.Range("W10").Value = 22
$ws.Range("Y10").Value = 45
$ws.Range("Z10").Value = 6.7
$ws.Range("AA10").Value = 6
$ws.Range("AB10").Value = 18
$ws.Range("AC10").Value = 110
$ws.Range("AI10").Value = 50

# Row 11 (Cheongju vs Asan)
$ws.Range("K11").Value = 9.5
$ws.Range("L11").Value = 1.3
$ws.Range("M11").Value = 3.4
$ws.Range("N11").Value = 2.03
$ws.Range("O11").Value = 1.78
$ws.Range("X11").Value = 26
$ws.Range("AD11").Value = 7.5
$ws.Range("AJ11").Value = 251

# Row 14 (Hartford Athletic vs North Carolina)
$ws.Range("H14").Value = 3.4
$ws.Range("K14").Value = 7.5
$ws.Range("X14").Value = 24
$ws.Range("Z14").Value = 7.5
